# Automatische test-sync: 2025-06-19 22:34:50
# Append a new mail-log row to the "Logs" sheet and bump the corresponding
# category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 47

$logs.Cells.Item($newRow, 1).Value = "Is product X op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 22:34:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 7

# Extend the conditional-formatting ranges (D and G columns) to cover the new row.
$catFcs = $logs.Range("D2:D46").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D47"))
}

$ansFcs = $logs.Range("G2:G46").FormatConditions
for ($i = 1; $i -le $ansFcs.Count; $i++) {
    $ansFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G47"))
}
